$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose data (id in column B, and columns F:AC) got swapped between them.
$pairs = @(@(180, 181), @(198, 199))

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    # Column B (id) swap
    $id1 = $ws.Cells.Item($row1, 2).Value2
    $id2 = $ws.Cells.Item($row2, 2).Value2
    $ws.Cells.Item($row1, 2).Value2 = $id2
    $ws.Cells.Item($row2, 2).Value2 = $id1

    # Columns F..AC (6..29) swap
    for ($col = 6; $col -le 29; $col++) {
        $v1 = $ws.Cells.Item($row1, $col).Value2
        $v2 = $ws.Cells.Item($row2, $col).Value2
        $ws.Cells.Item($row1, $col).Value2 = $v2
        $ws.Cells.Item($row2, $col).Value2 = $v1
    }
}
